$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LJ Speech")

# Update numeric values in column C
$ws.Range("C2").Value = 36
$ws.Range("C3").Value = 35
$ws.Range("C4").Value = 34
$ws.Range("C5").Value = 36
$ws.Range("C7").Value = 36
$ws.Range("C8").Value = 36
$ws.Range("C10").Value = 41
$ws.Range("C11").Value = 37
$ws.Range("C12").Value = 32
$ws.Range("C13").Value = 45
$ws.Range("C14").Value = 36
$ws.Range("C15").Value = 39
$ws.Range("C16").Value = 21

# Update inline string values in column B
$ws.Range("B8").Value = "<thow>"
$ws.Range("B12").Value = "<high>"
$ws.Range("B16").Value = "<tie>"
